$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A ("Match ID") and shift everything right.
$ws.Columns.Item(1).Insert()

# Header text for the new column (row 2 is the header row).
$ws.Range("A2").Value = "Match ID"

# Match the bold-header styling used by the rest of the header row.
$ws.Range("A2:A19").Font.Bold = $true

# Every player row (4-19) plus the hidden totals row (20) belongs to match 3.
$ws.Range("A4:A20").Value = 3

# Row 3 is a blank hidden spacer row; give it the same bold style as A2/A4:A19
# even though it has no value, matching the original authoring pattern.
$ws.Range("A3").Font.Bold = $true

# Re-fit the two hidden rows we touched so Excel drops the transient explicit
# row-height / customHeight it stamps on edited hidden rows.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(20).AutoFit()

# Restore the original selection anchor on the newly inserted column.
$null = $ws.Range("A2:A19").Select()
